$wb = $excel.ActiveWorkbook

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 225
$wsTrans.Range("A3").Value = 224
$wsTrans.Range("A4").Value = 222
$wsTrans.Range("A5").Value = 221
# Touch row 6 so it materializes as an empty (unstyled) row, extending the
# used range to A1:J6.
$wsTrans.Range("A6:J6").Style = "Normal"
$wsTrans.Range("J3").Select()

# --- Repayment Schedule sheet ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("C15").Select()

# --- Summary sheet (becomes the active tab) ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 3139.43
$wsSummary.Activate()
$wsSummary.Range("E4").Select()
